# Update hours (all team) - Time Sheet.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 updates
$ws.Range("D12").Value = 16
$ws.Range("F12").Value = 14

# Row 13 updates (new entries)
$ws.Range("B13").Value = 30
$ws.Range("D13").Value = 28
$ws.Range("F13").Value = 15

# Update the active cell / selection to B14 (as in the diff)
$ws.Range("B14").Select()
